$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like
# "23.385.26", "1.000", "15.50" keep their exact literal formatting
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "23.385.26"
$ws.Range("E2").Value = "  -1.38%  "
$ws.Range("D3").Value = "1.638.79"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "0.9991"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").Value = "298.98"
$ws.Range("E6").Value = "  -1.30%  "
$ws.Range("D7").Value = "0.3784"
$ws.Range("E7").Value = "  -1.56%  "
$ws.Range("D8").Value = "50.28"
$ws.Range("E8").Value = "  -1.64%  "
$ws.Range("D9").Value = "0.3483"
$ws.Range("E9").Value = "  -3.54%  "
$ws.Range("D10").Value = "0.08054"
$ws.Range("E10").Value = "  -1.79%  "
$ws.Range("D11").Value = "1.213"
$ws.Range("E11").Value = "  -1.43%  "
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").Value = "22.01"
$ws.Range("E13").Value = "  -1.72%  "
$ws.Range("D14").Value = "6.296"
$ws.Range("E14").Value = "  -2.40%  "
$ws.Range("D15").Value = "7.242"
$ws.Range("E15").Value = "  -2.73%  "
$ws.Range("D16").Value = "0.00001199"
$ws.Range("E16").Value = "  -1.83%  "
$ws.Range("D17").Value = "1.634.90"
$ws.Range("E17").Value = "  -1.30%  "
$ws.Range("D18").Value = "94.71"
$ws.Range("E18").Value = "  -2.99%  "
$ws.Range("D19").Value = "0.06960"
$ws.Range("E19").Value = "  -1.05%  "
$ws.Range("D20").Value = "6.609"
$ws.Range("E20").Value = "  -2.71%  "
$ws.Range("D21").Value = "17.30"
$ws.Range("E21").Value = "  -1.52%  "
$ws.Range("D22").Value = "0.9995"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").Value = "12.38"
$ws.Range("E23").Value = "  -2.72%  "
$ws.Range("D24").Value = "23.405.54"
$ws.Range("E24").Value = "  -1.33%  "
$ws.Range("D25").Value = "2.434"
$ws.Range("E25").Value = "  -2.27%  "
$ws.Range("D26").Value = "2.977"
$ws.Range("E26").Value = "  -1.63%  "
$ws.Range("D27").Value = "20.97"
$ws.Range("E27").Value = "  -1.40%  "
$ws.Range("D28").Value = "149.54"
$ws.Range("E28").Value = "  -2.68%  "
$ws.Range("D29").Value = "5.166"
$ws.Range("E29").Value = "  -1.48%  "
$ws.Range("D30").Value = "131.75"
$ws.Range("E30").Value = "  -1.68%  "
$ws.Range("D31").Value = "1.819.87"
$ws.Range("E31").Value = "  -1.00%  "
$ws.Range("D32").Value = "6.794"
$ws.Range("E32").Value = "  -4.61%  "
$ws.Range("D33").Value = "2.137"
$ws.Range("E33").Value = "  -5.15%  "
$ws.Range("D34").Value = "11.18"
$ws.Range("E34").Value = "  -7.03%  "
$ws.Range("D35").Value = "0.9890"
$ws.Range("E35").Value = "  -6.28%  "
$ws.Range("D36").Value = "0.02679"
$ws.Range("E36").Value = "  -4.58%  "
$ws.Range("D37").Value = "0.08756"
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("D38").Value = "0.2418"
$ws.Range("E38").Value = "  -3.37%  "
$ws.Range("D39").Value = "5.861"
$ws.Range("E39").Value = "  -3.72%  "
$ws.Range("D40").Value = "0.06806"
$ws.Range("E40").Value = "  -2.58%  "
$ws.Range("D41").Value = "12.76"
$ws.Range("E41").Value = "  -2.00%  "
$ws.Range("D42").Value = "0.6805"
$ws.Range("E42").Value = "  -2.64%  "
$ws.Range("D43").Value = "1.287"
$ws.Range("E43").Value = "  -3.70%  "
$ws.Range("D44").Value = "15.50"
$ws.Range("E44").Value = "  -2.83%  "
$ws.Range("D45").Value = "0.9981"
$ws.Range("E45").Value = "  -0.31%  "
$ws.Range("D46").Value = "0.6320"
$ws.Range("E46").Value = "  -2.89%  "
$ws.Range("D47").Value = "2.230"
$ws.Range("E47").Value = "  -2.86%  "
$ws.Range("D48").Value = "3.899"
$ws.Range("E48").Value = "  -1.64%  "

# Rows 49-50: Cronos/Quant swapped order with updated data
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "126.94"
$ws.Range("E49").Value = "  -0.95%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.07662"
$ws.Range("E50").Value = "  -2.96%  "

$ws.Range("D51").Value = "1.218"
$ws.Range("E51").Value = "  +2.00%  "

# Restore default (General) style on column D now that the exact
# text has been written, so no stray number-format style lingers.
$ws.Range("D2:D51").Style = "Normal"
